$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 38 inherits the same base formatting (borders, style) as row 37 directly above it.
$ws.Range("A37:E37").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122)

# Fill in the new watchlist-deletion test case.
$ws.Range("A38").Value = "TestCase_E37"
$ws.Range("C38").Value = "Verify that user is able to delete a watchlist||Verify that user is not able to see his watchlist on his own `nprofile page once that particular watchlist is deleted."
$ws.Range("B38").Value = "OPQA-625`n||OPQA-1104"
$ws.Range("D38").Value = "Y"
$ws.Range("E38").Value = "PASS"

# Jira id / Description columns wrap their (multi-line) text, same as other multi-line rows.
$ws.Range("B38").WrapText = $true
$ws.Range("C38").WrapText = $true

# Row grows to fit the two wrapped lines.
$ws.Rows.Item(38).RowHeight = 30

# Reflect the author's final cursor position on the sheet.
$ws.Range("B2").Select() | Out-Null
